$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '246.88'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '21.76'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.457'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05680'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.372'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8000'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.035'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'One'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0005891'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '8OneONE'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1449'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07241'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03163'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.02938'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09281'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001654'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '14BitForexTokenBF'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.206'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04711'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006367'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.005042'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0003200'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.800'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.431'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.125'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.3287'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04088'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006941'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.003500'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1038'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008970'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005831'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6825'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.01007'
